$wb = $excel.ActiveWorkbook

# --- Add the new trailing worksheet ("Sheet4") after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Populate row 1 with the new test-case data
$newSheet.Range("A1").Value = "selenium"
$newSheet.Range("B1").Value = "selenium@gmail.com"
$newSheet.Range("C1").Value = "apartment"
$newSheet.Range("D1").Value = "looking for apartment"
$newSheet.Range("E1").Value = 40000
$newSheet.Range("F1").Value = 2000
$newSheet.Range("G1").Value = 2
$newSheet.Range("H1").Value = 5

# Hyperlink the email address cell, then restore the shared Hyperlink cell style
$newSheet.Hyperlinks.Add($newSheet.Range("B1"), "mailto:selenium@gmail.com")
$newSheet.Range("B1").Style = "Hyperlink"

# Put the focus/selection on the last populated cell of the new sheet
$newSheet.Range("H1").Select()

# --- Switch the active tab back to UsersData ---
$usersData = $wb.Worksheets.Item("UsersData")
$usersData.Activate()
